$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: F1 "time_taken" (style copied from B1 header cell) ---
$ws.Range("F1").Value = 'time_taken'
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Rows 49-86: geneSymbol / geneName / geneConfidence reshuffled ---
$bVals4986 = @(
  'RNASEH2A',
  'RNASEH2B',
  'RNASEH2C',
  'SAMHD1',
  'SCN1A',
  'SCN8A',
  'SLC2A1',
  'SPAST',
  'SPG11',
  'SPTAN1',
  'TAF1',
  'TNR',
  'TREX1',
  'TUBA1A',
  'VPS13D',
  'ZC4H2',
  'ZSWIM6',
  'AGAP1',
  'ALK',
  'CACNA1A',
  'COL4A1',
  'DHX32',
  'GSX2',
  'HPDL',
  'NALCN',
  'NDUFA12',
  'PANK2',
  'PROC',
  'RHOB',
  'ASXL3',
  'ATP1A3',
  'BCL11A',
  'GAD1',
  'HPCA',
  'IQSEC2',
  'KANK1',
  'KDM7A',
  'KMT2B'
)
$cVals4986 = @(
  'ribonuclease H2 subunit A',
  'ribonuclease H2 subunit B',
  'ribonuclease H2 subunit C',
  'SAM and HD domain containing deoxynucleoside triphosphate triphosphohydrolase 1',
  'sodium voltage-gated channel alpha subunit 1',
  'sodium voltage-gated channel alpha subunit 8',
  'solute carrier family 2 member 1',
  'spastin',
  'SPG11, spatacsin vesicle trafficking associated',
  'spectrin alpha, non-erythrocytic 1',
  'TATA-box binding protein associated factor 1',
  'tenascin R',
  'three prime repair exonuclease 1',
  'tubulin alpha 1a',
  'vacuolar protein sorting 13 homolog D',
  'zinc finger C4H2-type containing',
  'zinc finger SWIM-type containing 6',
  'ArfGAP with GTPase domain, ankyrin repeat and PH domain 1',
  'ALK receptor tyrosine kinase',
  'calcium voltage-gated channel subunit alpha1 A',
  'collagen type IV alpha 1 chain',
  'DEAH-box helicase 32 (putative)',
  'GS homeobox 2',
  '4-hydroxyphenylpyruvate dioxygenase like',
  'sodium leak channel, non-selective',
  'NADH:ubiquinone oxidoreductase subunit A12',
  'pantothenate kinase 2',
  'protein C, inactivator of coagulation factors Va and VIIIa',
  'ras homolog family member B',
  'additional sex combs like 3, transcriptional regulator',
  'ATPase Na+/K+ transporting subunit alpha 3',
  'B-cell CLL/lymphoma 11A',
  'glutamate decarboxylase 1',
  'hippocalcin',
  'IQ motif and Sec7 domain 2',
  'KN motif and ankyrin repeat domains 1',
  'lysine demethylase 7A',
  'lysine methyltransferase 2B'
)
$dVals4986 = @(
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '3',
  '2',
  '2',
  '2',
  '2',
  '2',
  '2',
  '2',
  '2',
  '2',
  '2',
  '2',
  '2',
  '1',
  '1',
  '1',
  '1',
  '1',
  '1',
  '1',
  '1',
  '1'
)

for ($i = 0; $i -lt $bVals4986.Length; $i++) {
    $r = $i + 49
    $ws.Cells.Item($r, 2).Value = $bVals4986[$i]
    $ws.Cells.Item($r, 3).Value = $cVals4986[$i]
    $ws.Cells.Item($r, 4).Value = "'" + $dVals4986[$i]
}

# --- F column (time_taken) for all data rows 2-94 ---
$fVals = @(
  '2021-10-05 10:50:24.159513',
  '2021-10-05 10:50:24.159525',
  '2021-10-05 10:50:24.159528',
  '2021-10-05 10:50:24.159531',
  '2021-10-05 10:50:24.159533',
  '2021-10-05 10:50:24.159536',
  '2021-10-05 10:50:24.159539',
  '2021-10-05 10:50:24.159541',
  '2021-10-05 10:50:24.159544',
  '2021-10-05 10:50:24.159546',
  '2021-10-05 10:50:24.159549',
  '2021-10-05 10:50:24.159551',
  '2021-10-05 10:50:24.159554',
  '2021-10-05 10:50:24.159557',
  '2021-10-05 10:50:24.159559',
  '2021-10-05 10:50:24.159562',
  '2021-10-05 10:50:24.159564',
  '2021-10-05 10:50:24.159567',
  '2021-10-05 10:50:24.159569',
  '2021-10-05 10:50:24.159572',
  '2021-10-05 10:50:24.159574',
  '2021-10-05 10:50:24.159577',
  '2021-10-05 10:50:24.159579',
  '2021-10-05 10:50:24.159582',
  '2021-10-05 10:50:24.159584',
  '2021-10-05 10:50:24.159587',
  '2021-10-05 10:50:24.159589',
  '2021-10-05 10:50:24.159592',
  '2021-10-05 10:50:24.159594',
  '2021-10-05 10:50:24.159597',
  '2021-10-05 10:50:24.159599',
  '2021-10-05 10:50:24.159602',
  '2021-10-05 10:50:24.159605',
  '2021-10-05 10:50:24.159607',
  '2021-10-05 10:50:24.159610',
  '2021-10-05 10:50:24.159612',
  '2021-10-05 10:50:24.159615',
  '2021-10-05 10:50:24.159617',
  '2021-10-05 10:50:24.159619',
  '2021-10-05 10:50:24.159622',
  '2021-10-05 10:50:24.159625',
  '2021-10-05 10:50:24.159627',
  '2021-10-05 10:50:24.159630',
  '2021-10-05 10:50:24.159632',
  '2021-10-05 10:50:24.159634',
  '2021-10-05 10:50:24.159637',
  '2021-10-05 10:50:24.159640',
  '2021-10-05 10:50:24.159642',
  '2021-10-05 10:50:24.159645',
  '2021-10-05 10:50:24.159647',
  '2021-10-05 10:50:24.159650',
  '2021-10-05 10:50:24.159652',
  '2021-10-05 10:50:24.159655',
  '2021-10-05 10:50:24.159657',
  '2021-10-05 10:50:24.159660',
  '2021-10-05 10:50:24.159662',
  '2021-10-05 10:50:24.159665',
  '2021-10-05 10:50:24.159667',
  '2021-10-05 10:50:24.159670',
  '2021-10-05 10:50:24.159672',
  '2021-10-05 10:50:24.159675',
  '2021-10-05 10:50:24.159677',
  '2021-10-05 10:50:24.159680',
  '2021-10-05 10:50:24.159682',
  '2021-10-05 10:50:24.159686',
  '2021-10-05 10:50:24.159689',
  '2021-10-05 10:50:24.159691',
  '2021-10-05 10:50:24.159694',
  '2021-10-05 10:50:24.159696',
  '2021-10-05 10:50:24.159699',
  '2021-10-05 10:50:24.159701',
  '2021-10-05 10:50:24.159704',
  '2021-10-05 10:50:24.159706',
  '2021-10-05 10:50:24.159709',
  '2021-10-05 10:50:24.159711',
  '2021-10-05 10:50:24.159714',
  '2021-10-05 10:50:24.159718',
  '2021-10-05 10:50:24.159721',
  '2021-10-05 10:50:24.159723',
  '2021-10-05 10:50:24.159726',
  '2021-10-05 10:50:24.159728',
  '2021-10-05 10:50:24.159731',
  '2021-10-05 10:50:24.159733',
  '2021-10-05 10:50:24.159736',
  '2021-10-05 10:50:24.159738',
  '2021-10-05 10:50:24.159741',
  '2021-10-05 10:50:24.159744',
  '2021-10-05 10:50:24.159746',
  '2021-10-05 10:50:24.159749',
  '2021-10-05 10:50:24.159751',
  '2021-10-05 10:50:24.159754',
  '2021-10-05 10:50:24.159756',
  '2021-10-05 10:50:24.159760'
)
for ($i = 0; $i -lt $fVals.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 6).Value = $fVals[$i]
}

# --- New rows 87-94: append remaining gene rows ---
$aVals8794 = @(
  '85',
  '86',
  '87',
  '88',
  '89',
  '90',
  '91',
  '92'
)
$bVals8794 = @(
  'MAOB',
  'MFN2',
  'NEXMIF',
  'PAK3',
  'PCDH19',
  'SHANK3',
  'SMARCB1',
  'TENM1'
)
$cVals8794 = @(
  'monoamine oxidase B',
  'mitofusin 2',
  'neurite extension and migration factor',
  'p21 (RAC1) activated kinase 3',
  'protocadherin 19',
  'SH3 and multiple ankyrin repeat domains 3',
  'SWI/SNF related, matrix associated, actin dependent regulator of chromatin, subfamily b, member 1',
  'teneurin transmembrane protein 1'
)
$dVals8794 = @(
  '1',
  '1',
  '1',
  '1',
  '1',
  '1',
  '1',
  '1'
)
$eVals8794 = @(
  'Cerebral Palsy',
  'Cerebral Palsy',
  'Cerebral Palsy',
  'Cerebral Palsy',
  'Cerebral Palsy',
  'Cerebral Palsy',
  'Cerebral Palsy',
  'Cerebral Palsy'
)
$fVals8794 = @(
  '2021-10-05 10:50:24.159741',
  '2021-10-05 10:50:24.159744',
  '2021-10-05 10:50:24.159746',
  '2021-10-05 10:50:24.159749',
  '2021-10-05 10:50:24.159751',
  '2021-10-05 10:50:24.159754',
  '2021-10-05 10:50:24.159756',
  '2021-10-05 10:50:24.159760'
)

for ($i = 0; $i -lt $aVals8794.Length; $i++) {
    $r = $i + 87
    # Copy the full-row style template from the previous row (86) first so
    # borders/bold/alignment match the rest of the table.
    $destAddr = "A" + $r + ":E" + $r
    $ws.Range("A86:E86").Copy() | Out-Null
    $ws.Range($destAddr).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = "'" + $aVals8794[$i]
    $ws.Cells.Item($r, 2).Value = $bVals8794[$i]
    $ws.Cells.Item($r, 3).Value = $cVals8794[$i]
    $ws.Cells.Item($r, 4).Value = "'" + $dVals8794[$i]
    $ws.Cells.Item($r, 5).Value = $eVals8794[$i]
    $ws.Cells.Item($r, 6).Value = $fVals8794[$i]
}
$excel.CutCopyMode = 0

